# Backlog.xlsx update
#
# The backlog item that was in row 8 ("Implement load to working status")
# is reordered to the bottom of the "HIGH priority" block (new row 13);
# rows 9-13 each shift up by one row to fill the gap. While doing this
# reshuffle, the item that lands on the new row 8/9 ("Create HTML
# template to incorporate story" / "Make Play start game") gets its
# Status marked "Complete". Finally the selected cell moves to F15 and
# that (empty) cell picks up a freshly-applied format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# --- capture the "Neutral / thin grey border" look (currently on F11)
#     before we overwrite F11's contents below; it needs to land on the
#     new F10. ---
$ws.Range("F11").Copy()
$ws.Range("F10").PasteSpecial($xlPasteFormats)

# --- row 8 becomes what used to be row 9, with Status -> Complete ---
$ws.Range("A8").Value = "Create HTML template to incorporate story"
$ws.Range("B8").Value = "Task/Dev"
$ws.Range("C8").ClearContents()
$ws.Range("F7").Copy()
$ws.Range("F8").PasteSpecial($xlPasteFormats)
$ws.Range("F8").Value = "Complete"

# --- row 9 becomes what used to be row 10, with Status -> Complete ---
$ws.Range("A9").Value = "Make Play start game"
$ws.Range("B9").Value = "Use case"
$ws.Range("C9").Value = "Play button needs to start the game"
$ws.Range("F7").Copy()
$ws.Range("F9").PasteSpecial($xlPasteFormats)
$ws.Range("F9").Value = "Complete"

# --- row 10 becomes what used to be row 11, Status stays "In progress"
#     (format already copied onto F10 above) ---
$ws.Range("A10").Value = "Implement story in game"
$ws.Range("B10").Value = "Task/Dev"
$ws.Range("C10").Value = "Implement the story and the various of options"
$ws.Range("F10").Value = "In progress"

# --- row 11 becomes what used to be row 12, Status blank ---
$ws.Range("A11").Value = "Implement karma system"
$ws.Range("B11").Value = "Task/Dev"
$ws.Range("C11").Value = "Implemenet karma system to keep track of good/bad deeds"
$ws.Range("F14").Copy()
$ws.Range("F11").PasteSpecial($xlPasteFormats)
$ws.Range("F11").ClearContents()

# --- row 12 becomes what used to be row 13, Status blank ---
$ws.Range("A12").Value = "Implement save option in game"
$ws.Range("B12").Value = "Task/Dev"
$ws.Range("C12").Value = "Implement save option in game"
$ws.Range("F12").ClearContents()

# --- row 13 becomes the original row 8, Status blank ---
$ws.Range("A13").Value = "Implement load to working status"
$ws.Range("B13").Value = "Task/Dev"
$ws.Range("C13").Value = "Make load option load saved files"
$ws.Range("F13").ClearContents()

# --- F15 (still empty) picks up an explicit font/format touch ---
$ws.Range("F15").Font.ThemeColor = 1

# --- move the active selection to F15 ---
$ws.Range("F15").Select()
